# Auto-generated PowerShell COM-interop script applying the commit diff.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overall Completeness")
$data1 = @(
  @(100, 35202, 100, 10561, 100, 5427, 100, 6681, 100, 7484, 100, 5049),
  @(94.82, 33380, 94.71, 10002, 95.49, 5182, 95.84, 6403, 94.36, 7062, 93.7, 4731),
  @(100, 35202, 100, 10561, 100, 5427, 100, 6681, 100, 7484, 100, 5049),
  @(100, 35202, 100, 10561, 100, 5427, 100, 6681, 100, 7484, 100, 5049),
  @(100, 35202, 100, 10561, 100, 5427, 100, 6681, 100, 7484, 100, 5049),
  @(99.14, 34901, 98.7, 10424, 99.02, 5374, 99.27, 6632, 99.49, 7446, 99.52, 5025),
  @(100, 35202, 100, 10561, 100, 5427, 100, 6681, 100, 7484, 100, 5049),
  @(82.4, 29008, 65.21, 6887, 88.24, 4789, 89.79, 5999, 92.54, 6926, 87.28, 4407),
  @(29.88, 10518, 28.55, 3015, 37.44, 2032, 34.1, 2278, 27.9, 2088, 21.89, 1105),
  @(49.4, 17390, 48.24, 5095, 48.65, 2640, 53.09, 3547, 53.74, 4022, 41.32, 2086),
)
for ($r = 0; $r -lt $data1.Length; $r++) {
  $row = $data1[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws1.Cells.Item(2 + $r, 2 + $c).Value = $row[$c]
  }
}

$ws2 = $wb.Worksheets.Item("Breast variable Completeness")
$data2 = @(
  @(58.9, 7883, 31.7, 1751, 69.37, 1610, 84.85, 1820, 80.68, 1679, 77.91, 1023),
  @(34.45, 4611, 8.44, 466, 42.09, 977, 55.48, 1190, 56.46, 1175, 61.16, 803),
  @(48.64, 6510, 19.5, 1077, 55.58, 1290, 77.11, 1654, 75.25, 1566, 70.3, 923),
  @(3.24, 433, 0.02, 1, 3.83, 89, 7.79, 167, 6.97, 145, 2.36, 31),
  @(0.99, 132, 0, 0, 0.47, 11, 1.68, 36, 3.03, 63, 1.68, 22),
  @(59.25, 7929, 40.99, 2264, 65.83, 1528, 71, 1523, 74.58, 1552, 80.88, 1062),
  @(55.8, 7468, 36.79, 2032, 62.77, 1457, 66.95, 1436, 71.7, 1492, 80.05, 1051),
  @(99.13, 13266, 98.91, 5463, 99.22, 2303, 98.83, 2120, 99.57, 2072, 99.62, 1308),
  @(56.7, 7588, 37.15, 2052, 52.09, 1209, 73.15, 1569, 81.98, 1706, 80.12, 1052),
  @(53.52, 7163, 34.84, 1924, 50.11, 1163, 70.12, 1504, 76.65, 1595, 74.41, 977),
  @(40.27, 5389, 25.71, 1420, 37.18, 863, 49.04, 1052, 58.05, 1208, 64.43, 846),
  @(64.31, 8607, 50.15, 2770, 74.62, 1732, 80.75, 1732, 73.43, 1528, 64.36, 845),
  @(49.3, 6598, 36.61, 2022, 57.13, 1326, 61.21, 1313, 58.87, 1225, 54.23, 712),
  @(52.76, 7061, 39.65, 2190, 65.1, 1511, 64.38, 1381, 60.74, 1264, 54.46, 715),
)
for ($r = 0; $r -lt $data2.Length; $r++) {
  $row = $data2[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws2.Cells.Item(2 + $r, 2 + $c).Value = $row[$c]
  }
}

$ws3 = $wb.Worksheets.Item("Non-Null count by SITE (only 3)")
# Insert the 3 new rows (shifts existing rows down)
$ws3.Rows.Item(2).Insert()
$ws3.Rows.Item(2).ClearFormats()
$ws3.Rows.Item(4).Insert()
$ws3.Rows.Item(4).ClearFormats()
$ws3.Rows.Item(6).Insert()
$ws3.Rows.Item(6).ClearFormats()
$data3 = @(
  @("Bone and Articular Cartilage", 20, 3, 4),
  @("NMSC", 6681, 94, 347),
  @("Nasal cavity, middle ear, sinuses", 11, 10, 1),
  @("bladder", 121, 77, 3),
  @("bowel", 61, 41, 16),
  @("brain, CNS", 167, 0, 4),
  @("breast", 9950, 8205, 7904),
  @("cervix uteri", 175, 87, 51),
  @("colorectal cancer", 1398, 959, 171),
  @("corpus uteri", 707, 572, 120),
  @("gallbladder", 29, 21, 1),
  @("hodgkin lymphoma", 62, 44, 1),
  @("hypopharynx", 3, 2, 1),
  @("in-situ", 4567, 4241, 812),
  @("kaposi sarcoma", 0, 0, 0),
  @("kidney", 246, 153, 134),
  @("larynx", 5, 4, 1),
  @("leukemia", 255, 33, 0),
  @("lip, oral cavity", 112, 74, 52),
  @("liver", 73, 25, 15),
  @("lung", 582, 487, 234),
  @("melanoma of skin", 1176, 829, 324),
  @("mesothelioma", 28, 15, 0),
  @("multiple myeloma", 184, 30, 1),
  @("nasopharynx", 4, 3, 0),
  @("non-Hodgkin lymphoma", 484, 325, 2),
  @("oesophagus", 114, 64, 6),
  @("oropharynx", 36, 27, 10),
  @("ovary", 648, 505, 53),
  @("pancreas", 245, 157, 65),
  @("salivary glands", 26, 10, 12),
  @("stomach", 91, 54, 16),
  @("thyroid", 243, 80, 108),
  @("unknown/uncertain", 441, 131, 25),
  @("vagina", 12, 5, 1),
  @("vulva", 51, 23, 23),
)
for ($r = 0; $r -lt $data3.Length; $r++) {
  $row = $data3[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws3.Cells.Item(2 + $r, 1 + $c).Value = $row[$c]
  }
}

$ws4 = $wb.Worksheets.Item("Groups by SITE")
# Insert the 3 new rows (shifts existing rows down)
$ws4.Rows.Item(2).Insert()
$ws4.Rows.Item(2).ClearFormats()
$ws4.Rows.Item(4).Insert()
$ws4.Rows.Item(4).ClearFormats()
$ws4.Rows.Item(6).Insert()
$ws4.Rows.Item(6).ClearFormats()
$data4 = @(
  @("Bone and Articular Cartilage", 29),
  @("NMSC", 7779),
  @("Nasal cavity, middle ear, sinuses", 11),
  @("bladder", 145),
  @("bowel", 68),
  @("brain, CNS", 198),
  @("breast", 11590),
  @("cervix uteri", 281),
  @("colorectal cancer", 1616),
  @("corpus uteri", 806),
  @("gallbladder", 33),
  @("hodgkin lymphoma", 131),
  @("hypopharynx", 4),
  @("in-situ", 6184),
  @("kaposi sarcoma", 1),
  @("kidney", 303),
  @("larynx", 5),
  @("leukemia", 314),
  @("lip, oral cavity", 124),
  @("liver", 85),
  @("lung", 661),
  @("melanoma of skin", 1611),
  @("mesothelioma", 31),
  @("multiple myeloma", 213),
  @("nasopharynx", 5),
  @("non-Hodgkin lymphoma", 583),
  @("oesophagus", 129),
  @("oropharynx", 38),
  @("ovary", 801),
  @("pancreas", 274),
  @("salivary glands", 42),
  @("stomach", 104),
  @("thyroid", 317),
  @("unknown/uncertain", 604),
  @("vagina", 18),
  @("vulva", 64),
)
for ($r = 0; $r -lt $data4.Length; $r++) {
  $row = $data4[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws4.Cells.Item(2 + $r, 1 + $c).Value = $row[$c]
  }
}

$ws5 = $wb.Worksheets.Item("Groups by Source")
$data5 = @(
  @(27332),
  @(202),
  @(2547),
  @(208),
  @(1),
  @(242),
  @(3700),
  @(138),
  @(13),
  @(147),
  @(15),
  @(334),
  @(149),
  @(7),
  @(1),
  @(2),
  @(13),
  @(70),
  @(27),
  @(54),
)
for ($r = 0; $r -lt $data5.Length; $r++) {
  $row = $data5[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws5.Cells.Item(2 + $r, 2 + $c).Value = $row[$c]
  }
}

$ws6 = $wb.Worksheets.Item("Groups by Source (Breast)")
$data6 = @(
  @(10072),
  @(81),
  @(972),
  @(63),
  @(1),
  @(222),
  @(1345),
  @(137),
  @(13),
  @(334),
  @(7),
  @(1),
  @(12),
  @(59),
  @(12),
  @(52),
)
for ($r = 0; $r -lt $data6.Length; $r++) {
  $row = $data6[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws6.Cells.Item(2 + $r, 2 + $c).Value = $row[$c]
  }
}

$ws7 = $wb.Worksheets.Item("YearOfDiag by Source (Breast)")
$data7 = @(
  @(83.32, 81.9, 85.54, 87.44, 63.23),
  @(0.96, 0.52, 0.73, 0, 0.62),
  @(0.08, 0.32, 0.32, 0, 17.3),
  @(0.08, 0.1, 0, 0, 1.17),
  @(0, 0.04, 0, 0, 0),
  @(1.33, 1.25, 0.98, 0, 1.4),
  @(9.89, 9.54, 7.11, 7.42, 12.54),
  @(0.17, 0.42, 0, 0, 2.25),
  @(0.21, 0.13, 0, 0, 0.08),
  @(2.62, 4.89, 4.63, 5.14, 0.05),
  @(0.06, 0, 0, 0, 0.09),
  @(0, 0, 0, 0, 0.01),
  @(0.05, 0.19, 0.26, 0, 0.01),
  @(0.7, 0.42, 0.13, 0, 0.56),
  @(0.05, 0.29, 0.29, 0, 0),
  @(0.49, 0, 0, 0, 0.7),
)
for ($r = 0; $r -lt $data7.Length; $r++) {
  $row = $data7[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws7.Cells.Item(2 + $r, 2 + $c).Value = $row[$c]
  }
}

